$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with pristine default style, used to restore styling
# after forcing text format on numeric-looking Price values.
$refStyle = $ws.Range("D2").Style

$ws.Range("D2").Value = '74.887.76'
$ws.Range("E2").Value = '  +1.50%  '

$ws.Range("D3").Value = '2.797.55'
$ws.Range("E3").Value = '  +7.03%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = $refStyle
$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '187.86'
$ws.Range("D5").Style = $refStyle
$ws.Range("E5").Value = '  +1.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '590.34'
$ws.Range("D6").Style = $refStyle
$ws.Range("E6").Value = '  +1.66%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = $refStyle
$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.544'
$ws.Range("D8").Style = $refStyle
$ws.Range("E8").Value = '  +3.16%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.189'
$ws.Range("D9").Style = $refStyle
$ws.Range("E9").Value = '  -2.05%  '

$ws.Range("D10").Value = '2.793.52'
$ws.Range("E10").Value = '  +6.87%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.374'
$ws.Range("D11").Style = $refStyle
$ws.Range("E11").Value = '  +6.08%  '

$ws.Range("E12").Value = '  -2.09%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.84'
$ws.Range("D13").Style = $refStyle
$ws.Range("E13").Value = '  +3.92%  '

$ws.Range("D14").Value = '3.301.53'
$ws.Range("E14").Value = '  +6.62%  '

$ws.Range("D15").Value = '74.826.93'
$ws.Range("E15").Value = '  +1.36%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000186'
$ws.Range("D16").Style = $refStyle
$ws.Range("E16").Value = '  +0.42%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.77'
$ws.Range("D17").Style = $refStyle
$ws.Range("E17").Value = '  +2.59%  '

$ws.Range("D18").Value = '2.789.78'
$ws.Range("E18").Value = '  +6.64%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.96'
$ws.Range("D19").Style = $refStyle
$ws.Range("E19").Value = '  -1.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.23'
$ws.Range("D20").Style = $refStyle
$ws.Range("E20").Value = '  +4.62%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '375.53'
$ws.Range("D21").Style = $refStyle
$ws.Range("E21").Value = '  +3.78%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.28'
$ws.Range("D22").Style = $refStyle
$ws.Range("E22").Value = '  +0.49%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.07'
$ws.Range("D23").Style = $refStyle
$ws.Range("E23").Value = '  +0.31%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("D24").Style = $refStyle
$ws.Range("E24").Value = '  +0.30%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.70'
$ws.Range("D25").Style = $refStyle
$ws.Range("E25").Value = '  +2.02%  '

$ws.Range("D26").Value = '2.946.27'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.13'
$ws.Range("D27").Style = $refStyle
$ws.Range("E27").Value = '  +1.76%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.60'
$ws.Range("D28").Style = $refStyle
$ws.Range("E28").Value = '  +4.26%  '

$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000102'
$ws.Range("D29").Style = $refStyle
$ws.Range("E29").Value = '  +10.83%  '

$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = $refStyle
$ws.Range("E30").Value = '  -0.67%  '

$ws.Range("E31").Value = '  +1.03%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '508.52'
$ws.Range("D32").Style = $refStyle
$ws.Range("E32").Value = '  -2.53%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.59'
$ws.Range("D33").Style = $refStyle
$ws.Range("E33").Value = '  +0.42%  '

$ws.Range("E34").Value = '  +3.46%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("D35").Style = $refStyle
$ws.Range("E35").Value = '  -0.09%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '165.05'
$ws.Range("D36").Style = $refStyle
$ws.Range("E36").Value = '  +1.82%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.72'
$ws.Range("D37").Style = $refStyle
$ws.Range("E37").Value = '  +4.18%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.117'
$ws.Range("D38").Style = $refStyle
$ws.Range("E38").Value = '  -0.26%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.36'
$ws.Range("D39").Style = $refStyle
$ws.Range("E39").Value = '  +0.49%  '

$ws.Range("E40").Value = '  -0.08%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '180.17'
$ws.Range("D41").Style = $refStyle
$ws.Range("E41").Value = '  +11.94%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.340'
$ws.Range("D42").Style = $refStyle
$ws.Range("E42").Value = '  +5.40%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.96'
$ws.Range("D43").Style = $refStyle
$ws.Range("E43").Value = '  +2.40%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.65'
$ws.Range("D44").Style = $refStyle
$ws.Range("E44").Value = '  +0.52%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '40.11'
$ws.Range("D45").Style = $refStyle
$ws.Range("E45").Value = '  +3.21%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.20'
$ws.Range("D46").Style = $refStyle
$ws.Range("E46").Value = '  +3.36%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0863'
$ws.Range("D47").Style = $refStyle
$ws.Range("E47").Value = '  +2.45%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.32'
$ws.Range("D48").Style = $refStyle
$ws.Range("E48").Value = '  -0.48%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.567'
$ws.Range("D49").Style = $refStyle
$ws.Range("E49").Value = '  +9.63%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.70'
$ws.Range("D50").Style = $refStyle
$ws.Range("E50").Value = '  +3.60%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.634'
$ws.Range("D51").Style = $refStyle
$ws.Range("E51").Value = '  +8.58%  '
